$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell while keeping it
# stored as TEXT (matching the source data, which is plain text like
# "496.03", not a real number). Briefly forcing a Text number format
# stops Excel from auto-converting the assigned string to a number;
# the format is then restored so the cell keeps its original (default)
# style, only the cached value type changes to string.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "54.119.85"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "2.259.53"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  +0.19%  "
Set-TextValue $ws.Range("D5") "496.03"
$ws.Range("E5").Value = "  +0.07%  "
Set-TextValue $ws.Range("D6") "128.42"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("E10").Value = "  +0.85%  "
Set-TextValue $ws.Range("D11") "0.335"
$ws.Range("E11").Value = "  +3.01%  "
Set-TextValue $ws.Range("D12") "4.76"
$ws.Range("E12").Value = "  +2.88%  "
$ws.Range("D13").Value = "2.657.84"
$ws.Range("E13").Value = "  -0.56%  "
Set-TextValue $ws.Range("D14") "22.72"
$ws.Range("E14").Value = "  +4.37%  "
$ws.Range("D15").Value = "54.103.80"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "2.265.37"
$ws.Range("E17").Value = "  -1.17%  "
Set-TextValue $ws.Range("D18") "10.23"
$ws.Range("E18").Value = "  +1.66%  "
Set-TextValue $ws.Range("D19") "4.13"
$ws.Range("E19").Value = "  +1.00%  "
Set-TextValue $ws.Range("D20") "301.30"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E21").Value = "  -2.52%  "
Set-TextValue $ws.Range("D22") "0.999"
$ws.Range("E22").Value = "  +0.15%  "
Set-TextValue $ws.Range("D23") "60.76"
$ws.Range("E23").Value = "  -2.50%  "
Set-TextValue $ws.Range("D24") "0.996"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  -1.45%  "
Set-TextValue $ws.Range("D26") "7.29"
$ws.Range("E26").Value = "  +2.74%  "
Set-TextValue $ws.Range("D27") "170.68"
$ws.Range("E27").Value = "  +1.21%  "
Set-TextValue $ws.Range("D28") "1.61"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +0.13%  "
Set-TextValue $ws.Range("D33") "17.70"
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  -0.23%  "
Set-TextValue $ws.Range("D35") "0.940"
$ws.Range("E35").Value = "  +8.23%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  -0.68%  "
Set-TextValue $ws.Range("D40") "3.36"
$ws.Range("E40").Value = "  +0.38%  "
Set-TextValue $ws.Range("D41") "125.28"
$ws.Range("E41").Value = "  -2.53%  "
Set-TextValue $ws.Range("D42") "4.79"
$ws.Range("E42").Value = "  -0.03%  "
Set-TextValue $ws.Range("D43") "0.0492"
$ws.Range("E43").Value = "  +1.47%  "
Set-TextValue $ws.Range("D44") "0.0889"
$ws.Range("E44").Value = "  +0.02%  "
Set-TextValue $ws.Range("D45") "0.544"
$ws.Range("E45").Value = "  +0.15%  "
Set-TextValue $ws.Range("D46") "240.67"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("E49").Value = "  +0.30%  "
Set-TextValue $ws.Range("D50") "16.13"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("E51").Value = "  -0.83%  "
